# Daily update at 8 AM UTC
# Adds the new day's row (row 82) to the Wins Over Time tracking sheet,
# and moves the "latest row" date-format highlight from row 81 to row 82.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 81 is no longer the last row, so it reverts to the regular
# "YYYY-MM-DD HH:MM:SS" date number format used by every other data row.
$ws.Range("A81").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 82.
$ws.Range("A82").Value = 45822
$ws.Range("B82").Value = 349
$ws.Range("C82").Value = 350
$ws.Range("D82").Value = 356

# Row 82 is now the last row, so it gets the distinct "YYYY-MM-DD"
# date number format previously used by row 81.
$ws.Range("A82").NumberFormat = "YYYY-MM-DD"
